$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = -12.51290000000001
$ws.Range("D5").Value = -7.447500000000007
$ws.Range("B8").Value = 5.583599999999997
$ws.Range("D8").Value = -8.041300000000005
$ws.Range("B10").Value = 6.6314
$ws.Range("C11").Value = -14.0056
$ws.Range("B12").Value = 5.639799999999999
$ws.Range("C12").Value = -13.3249
$ws.Range("D12").Value = -7.915100000000003
$ws.Range("D13").Value = -8.595499999999994
$ws.Range("C15").Value = -14.41659999999999
$ws.Range("D15").Value = -8.359200000000001
$ws.Range("C17").Value = -14.0948
$ws.Range("B18").Value = 6.482699999999997
$ws.Range("D21").Value = -8.145299999999992
$ws.Range("B25").Value = 6.040899999999996
$ws.Range("D25").Value = -7.3426
$ws.Range("C26").Value = -12.6253
$ws.Range("C27").Value = -13.4216
$ws.Range("C28").Value = -13.8163
$ws.Range("C32").Value = -13.3129
$ws.Range("D32").Value = -9.004300000000001
$ws.Range("D36").Value = -7.880500000000001
$ws.Range("B37").Value = 8.6555
$ws.Range("C37").Value = -12.85019999999999
$ws.Range("D38").Value = -7.917199999999998
$ws.Range("C41").Value = -12.84870000000001
$ws.Range("D41").Value = -8.1251
$ws.Range("C47").Value = -12.4821
$ws.Range("D50").Value = -8.327499999999999
$ws.Range("C51").Value = -11.8021
$ws.Range("D52").Value = -7.8508
$ws.Range("B55").Value = 5.946499999999997
$ws.Range("D59").Value = -8.433299999999997
$ws.Range("C65").Value = -12.3999
$ws.Range("D67").Value = -6.998999999999994
$ws.Range("B68").Value = 4.740699999999995
$ws.Range("C73").Value = -11.51990000000001
$ws.Range("B77").Value = 8.717900000000002
$ws.Range("B78").Value = 9.173999999999994
$ws.Range("B79").Value = 8.991400000000004
$ws.Range("B80").Value = 9.284099999999997
$ws.Range("B81").Value = 6.093300000000002
$ws.Range("B82").Value = 5.779599999999997
$ws.Range("B84").Value = 6.6459
$ws.Range("C84").Value = -12.18440000000001
$ws.Range("D84").Value = -8.139199999999997
$ws.Range("C85").Value = -13.2659
$ws.Range("D86").Value = -8.465099999999998
$ws.Range("D88").Value = -7.835899999999999
$ws.Range("C89").Value = -14.1623
$ws.Range("D89").Value = -8.098499999999996
$ws.Range("C93").Value = -10.2265
$ws.Range("C95").Value = -13.21369999999999
$ws.Range("D95").Value = -7.859600000000004
$ws.Range("C98").Value = -13.3251
$ws.Range("C99").Value = -12.2884
$ws.Range("B101").Value = 5.595099999999997
$ws.Range("C101").Value = -13.75340000000001
$ws.Range("B102").Value = 6.946099999999998
$ws.Range("C102").Value = -12.36860000000001
$ws.Range("D105").Value = -7.203899999999996
